$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: locate the (first) paragraph whose text contains $needle.
# ------------------------------------------------------------------
function Find-ParagraphContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) "Исх № Т{num} от 08.08.2017 г."  ->  "Исх № Т{num} от {date} г."
# ------------------------------------------------------------------
$pNum = Find-ParagraphContaining $d "Исх"
$rNum = $pNum.Range
$fNum = $rNum.Find
$fNum.ClearFormatting()
$fNum.Text = "08.08.2017"
$fNum.Forward = $true
$fNum.Wrap = 0
if ($fNum.Execute()) {
    $rNum.Text = "{date}"
}

# ------------------------------------------------------------------
# 2) "Дата записи: 08.08.2017 г. "  ->  "Дата записи: {date} г. "
#    The hidden "_GoBack" bookmark used to sit around "{file}" in the
#    "Наименование: {file}," line; after the edit it instead wraps the
#    new "{date}" placeholder here.
# ------------------------------------------------------------------
$pDate = Find-ParagraphContaining $d "Дата записи"
$rDate = $pDate.Range
$fDate = $rDate.Find
$fDate.ClearFormatting()
$fDate.Text = "08.08.2017"
$fDate.Forward = $true
$fDate.Wrap = 0
if ($fDate.Execute()) {
    $rDate.Text = "{date}"
    $newBmStart = $rDate.Start
    $newBmEnd = $rDate.End

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $bmRange = $d.Range($newBmStart, $newBmEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
